# Insert a new row at position 361 (shifts existing rows 361..485 down to 362..486)
# and populate it with a new weekly price record for Espinaca / Femacal de La Calera.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(361).Insert()

$ws.Cells.Item(361, 1).Value = 3
$ws.Cells.Item(361, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(361, 3).Value = "Coquimbo"
$ws.Cells.Item(361, 4).Value = 44988
$ws.Cells.Item(361, 5).Value = 5
$ws.Cells.Item(361, 6).Value = 100112012
$ws.Cells.Item(361, 7).Value = "Espinaca"
$ws.Cells.Item(361, 8).Value = "Sin especificar"
$ws.Cells.Item(361, 9).Value = "Primera"
$ws.Cells.Item(361, 10).Value = 230
$ws.Cells.Item(361, 11).Value = 6000
$ws.Cells.Item(361, 12).Value = 6500
$ws.Cells.Item(361, 13).Value = 6239
$ws.Cells.Item(361, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(361, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(361, 16).Value = 2080
$ws.Cells.Item(361, 17).Value = 3
$ws.Cells.Item(361, 18).Value = "Hortaliza"
